$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-33 from 45208 to 45212
$ws.Range("C2:C33").Value = 45212
